# "Updated running SDLT CODE with bacs submit date"
# Rows 3 and 4 have been fully processed by the bot (bacs submitted), so their
# case-reference / description text and their tracked-step checkboxes are
# cleared back out, leaving just the coloured placeholder cells behind.
# Row 5 is advanced to the next SDLT case code waiting to be processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: clear the case code / description and the D:K tracking checkboxes.
$ws.Range("A3:B3").ClearContents()
$ws.Range("D3:K3").ClearContents()

# Row 4: clear the case code / description and the D:K tracking checkboxes.
$ws.Range("A4:B4").ClearContents()
$ws.Range("D4:K4").ClearContents()

# Row 5: move on to the next running SDLT case code (description / status
# checkboxes for this row are left as they were).
$ws.Range("A5").Value = "NBT1893"
